$wb = $excel.ActiveWorkbook

# --- 1. Duplicate the "2022-Q3" sheet to create the new "2022-Q4" sheet,
#        inserted right before it (so the tab order becomes:
#        总计, 2022-Q4, 2022-Q3, 2022-Q2) ---
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# --- 2. Update the new "2022-Q4" sheet with the latest-quarter figures.
#        Columns D:G hold numeric-looking values that are stored as TEXT
#        in this workbook, so force a text number format before writing
#        and reset the style afterwards to avoid leaving stray formatting. ---
$q4.Range("C2").Value = "华宝英国富时100指数A"
foreach ($addr in @("D2", "E2", "F2", "G2")) {
    $q4.Range($addr).NumberFormat = "@"
}
$q4.Range("D2").Value = "0.14"
$q4.Range("E2").Value = "94.75"
$q4.Range("F2").Value = "2.86"
$q4.Range("G2").Value = "0.0040"
foreach ($addr in @("D2", "E2", "F2", "G2")) {
    $q4.Range($addr).Style = "Normal"
}

$q4.Range("C3").Value = "华宝英国富时100指数C"
foreach ($addr in @("E3", "F3", "G3")) {
    $q4.Range($addr).NumberFormat = "@"
}
$q4.Range("E3").Value = "94.75"
$q4.Range("F3").Value = "2.86"
$q4.Range("G3").Value = "0.0023"
foreach ($addr in @("E3", "F3", "G3")) {
    $q4.Range($addr).Style = "Normal"
}

# --- 3. Update the "总计" (total) sheet: shift the existing two rows down
#        one slot and insert the new Q4 entry at the top of the data. ---
$total = $wb.Worksheets.Item("总计")
$total.Range("B3").Value = "2022-Q3"
$total.Range("B2").Value = "2022-Q4"

$total.Range("A3:D3").Copy()
$total.Range("A4:D4").PasteSpecial(-4122)
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.01

# --- 4. Copying "2022-Q3" made the new sheet the active tab; restore the
#        original selected tab ("2022-Q2"), which this edit leaves untouched. ---
$wb.Worksheets.Item("2022-Q2").Activate()
